$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.730.23'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.465.80'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '611.19'
$ws.Range("E5").Value = '  +1.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.64'
$ws.Range("E6").Value = '  -1.85%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.462.66'
$ws.Range("E7").Value = '  -0.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.595'
$ws.Range("E8").Value = '  -1.93%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.192'
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.12'
$ws.Range("E11").Value = '  -1.92%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.565'
$ws.Range("E12").Value = '  -2.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '44.50'
$ws.Range("E13").Value = '  -2.83%  '
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.020.55'
$ws.Range("E15").Value = '  -1.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.21'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.465.57'
$ws.Range("E17").Value = '  -1.26%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '585.02'
$ws.Range("E18").Value = '  -2.67%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.763.82'
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.853'
$ws.Range("E22").Value = '  -1.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.00'
$ws.Range("E23").Value = '  -1.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '95.93'
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '15.26'
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("E26").Value = '  -1.50%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.46'
$ws.Range("E28").Value = '  -3.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.02'
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.66'
$ws.Range("E30").Value = '  -2.78%  '
$ws.Range("E31").Value = '  -1.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.83'
$ws.Range("E32").Value = '  -5.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.24'
$ws.Range("E33").Value = '  -2.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.59'
$ws.Range("E34").Value = '  -3.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '587.46'
$ws.Range("E35").Value = '  -16.46%  '
$ws.Range("B36").Value = 'Cosmos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.62'
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0479'
$ws.Range("E37").Value = '  +2.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0962'
$ws.Range("E38").Value = '  -2.89%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '56.28'
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.140'
$ws.Range("E41").Value = '  -0.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.18'
$ws.Range("E42").Value = '  -9.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.249.58'
$ws.Range("E43").Value = '  -2.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0700'
$ws.Range("E44").Value = '  +2.58%  '
$ws.Range("E45").Value = '  -4.81%  '
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.80'
$ws.Range("E46").Value = '  -3.28%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '30.89'
$ws.Range("E47").Value = '  -3.54%  '
$ws.Range("E48").Value = '  -4.28%  '
$ws.Range("E49").Value = '  -2.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '134.03'
$ws.Range("E50").Value = '  +1.42%  '